$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Merge split runs "Number of outcomes for die " + "B" + " =6"
#    into a single run (visible text is unchanged, only the backing
#    run structure is simplified away).
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute("Number of outcomes for die B =6", $true, $false, $false, $false, $false, $true, 1, $false, "Number of outcomes for die B =6", 2)

# ------------------------------------------------------------------
# 2) Merge split runs "QUESTION " + "3." into a single run.
# ------------------------------------------------------------------
$null = $d.Content.Find.Execute("QUESTION 3.", $true, $false, $false, $false, $false, $true, 1, $false, "QUESTION 3.", 2)

# ------------------------------------------------------------------
# 3) The paragraph right after the page break currently holds just a
#    lone tab character (pPr has the 6624 tab stop; run is
#    <w:lastRenderedPageBreak/><w:tab/>). We want that run to instead
#    read <w:lastRenderedPageBreak/><w:t>PART-B</w:t>. Editing the
#    run in place drops the lastRenderedPageBreak marker, so instead
#    we insert a brand-new paragraph (carrying its own
#    lastRenderedPageBreak + "PART-B" text) right before it, then
#    delete the old tab-only paragraph outright.
# ------------------------------------------------------------------
$tabPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $ctext = $cand.Range.Text
    if ($ctext.Length -eq 2 -and $ctext[0] -eq "`t") {
        $tabPara = $cand
        break
    }
}
if ($tabPara -eq $null) {
    throw "could not locate the lastRenderedPageBreak/tab paragraph"
}

$insertionRange = $tabPara.Previous(1).Range

$newPartB = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:tabs><w:tab w:val="left" w:pos="6624"/></w:tabs></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>PART-B</w:t></w:r></w:p>
'@
$insertionRange.InsertXML($newPartB)

$tabPara.Range.Delete()

# ------------------------------------------------------------------
# 4) Insert the large new block of paragraphs (LOGIC-, Die A / Die B
#    descriptions, the write-up paragraphs, MATHEMATICAL CALCULATION-
#    and the P(Sum = n) = x/36 lines) right after "PART-B", and
#    remove the old now-superseded empty "tabs" paragraph that used
#    to directly follow it.
# ------------------------------------------------------------------
$partB = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd() -eq "PART-B") {
        $partB = $cand
        break
    }
}
if ($partB -eq $null) {
    throw "could not locate the freshly inserted PART-B paragraph"
}

$oldEmptyTabsPara = $partB.Next(1)

$bigFragment = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:tabs><w:tab w:val="left" w:pos="6624"/></w:tabs></w:pPr><w:r><w:t>LOGIC-</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Die A: Represents a standard six-sided die [1, 2, 3, 4, 5, 6].</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Die B: Also represents a standard six-sided die [1, 2, 3, 4, 5, 6].</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>The sum 2 has a probability of 1/36, while the sums 3, 4, 5, 6, and so on, increase incrementally up to 1/36 for sum 12</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Die A has a limitation where no face can have more than 4 spots.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Calculate the original probabilities of obtaining each sum from rolling Die A and Die B.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Generate the current counts of each sum based on the given dice.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Determine the target counts for each sum based on the original probabilities.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Iterate through Die A's faces, checking for faces with more than 4 spots.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:r><w:tab/></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>MATHEMATICAL CALCULATION-</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>P(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Sum = 2) = 1/36</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>P(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Sum = 3) = 2/36</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>P(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Sum = 4) = 3/36</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>P(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Sum = 5) = 4/36</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>P(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Sum = 6) = 5/36</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>P(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Sum = 7) = 6/36</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>P(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Sum = 8) = 5/36</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>P(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Sum = 9) = 4/36</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>P(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Sum = 10) = 3/36</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>P(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Sum = 11) = 2/36</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>P(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:eastAsia="en-IN"/></w:rPr><w:t>Sum = 12) = 1/36</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
'@
$partB.Range.InsertXML($bigFragment)

$oldEmptyTabsPara.Range.Delete()

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
